$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are textual (dotted thousand separators, leading
# zeros, etc.) in the source data, so force text format before assigning to
# prevent Excel from auto-coercing them to numbers; then restore the default
# "Normal" style so no stray per-cell formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.841.70'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.773.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.12%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +1.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '339.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.83%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3818'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3404'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.36%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.33'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.137'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.85%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07373'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.44'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.005'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.341'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.396'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.774.79'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001071'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06662'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.45%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.09'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.57%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.71%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.423'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.861.69'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.03'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.95%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.398'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.490'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.72%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.398'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.91%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '152.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.968.18'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.78%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '133.53'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.039'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.013'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.49%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08871'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.63'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.68%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02387'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.85%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6800'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.52%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06360'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.56%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.273'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2159'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.498'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.225'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.50%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.171'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.13%  '

$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.003'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.86%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6213'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.84%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.860'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.13%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.12'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.059'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.83%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07366'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.92%  '

$ws.Range("E51").Value = '  +4.07%  '
